$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit -------------------------------------------------------
# The "Case/Study/Demographic/Diagnosis" Cypher query (cell B2) dropped its
# trailing OPTIONAL MATCH (co:cohort) / `Cohort` output column.
$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and demo.breed in ['Airedale Terrier', 'Labrador Retriever','West Highland White Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder, Urethra'] and diag.best_response in ['Not Determined']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $newQuery

# --- View state -----------------------------------------------------------
# Scroll the sheet up one row and move the selection from D4 to B2, matching
# where the author was working after trimming the query text.
$ws.Activate()
$ws.Range("B2").Select()
